$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 248
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 14
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 2231
$ws.Range("L2").Value = 1326
$ws.Range("M2").Value = 905
$ws.Range("N2").Value = 900
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 55
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = -6
$ws.Range("S2").Value = -30
$ws.Range("T2").Value = 6
$ws.Range("U2").Value = 1
$ws.Range("V2").Value = 893
$ws.Range("W2").Value = 5.22
$ws.Range("X2").Value = 3.86
$ws.Range("Y2").Value = $null
$ws.Range("Z2").Value = $null
$ws.Range("AA2").Value = 146.58
$ws.Range("AB2").Value = 1532.56
$ws.Range("AC2").Value = 98
$ws.Range("AD2").Value = $null
$ws.Range("AE2").Value = 8152
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = $null
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 11041708
# Row 3
$ws.Range("D3").Value = 2485
$ws.Range("E3").Value = 101
$ws.Range("F3").Value = 101
$ws.Range("G3").Value = 35
$ws.Range("H3").Value = 15
$ws.Range("I3").Value = 17
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 2224
$ws.Range("L3").Value = 1298
$ws.Range("M3").Value = 926
$ws.Range("N3").Value = 926
$ws.Range("O3").Value = $null
$ws.Range("P3").Value = 55
$ws.Range("Q3").Value = 295
$ws.Range("R3").Value = -119
$ws.Range("S3").Value = 123
$ws.Range("T3").Value = 123
$ws.Range("U3").Value = 172
$ws.Range("V3").Value = 1046
$ws.Range("W3").Value = 4.06
$ws.Range("X3").Value = 0.59
$ws.Range("Y3").Value = 1.84
$ws.Range("Z3").Value = 0.65
$ws.Range("AA3").Value = 140.22
$ws.Range("AB3").Value = 1570.63
$ws.Range("AC3").Value = 152
$ws.Range("AD3").Value = 40.62
$ws.Range("AE3").Value = 8388
$ws.Range("AF3").Value = 0.74
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 1.62
$ws.Range("AI3").Value = 65.59
$ws.Range("AJ3").Value = 11041708
# Row 4
$ws.Range("D4").Value = 2291
$ws.Range("E4").Value = 140
$ws.Range("F4").Value = 140
$ws.Range("G4").Value = 119
$ws.Range("H4").Value = 96
$ws.Range("I4").Value = 96
$ws.Range("J4").Value = $null
$ws.Range("K4").Value = 2282
$ws.Range("L4").Value = 1277
$ws.Range("M4").Value = 1005
$ws.Range("N4").Value = 1005
$ws.Range("O4").Value = $null
$ws.Range("P4").Value = 55
$ws.Range("Q4").Value = 324
$ws.Range("R4").Value = -119
$ws.Range("S4").Value = -225
$ws.Range("T4").Value = 115
$ws.Range("U4").Value = 209
$ws.Range("V4").Value = 819
$ws.Range("W4").Value = 6.11
$ws.Range("X4").Value = 4.19
$ws.Range("Y4").Value = 9.949999999999999
$ws.Range("Z4").Value = 4.26
$ws.Range("AA4").Value = 127.02
$ws.Range("AB4").Value = 1714.18
$ws.Range("AC4").Value = 870
$ws.Range("AD4").Value = 10.67
$ws.Range("AE4").Value = 9107
$ws.Range("AF4").Value = 1.02
$ws.Range("AG4").Value = 120
$ws.Range("AH4").Value = 1.29
$ws.Range("AI4").Value = 13.79
$ws.Range("AJ4").Value = 11041708
# Row 5
$ws.Range("D5").Value = 3245
$ws.Range("E5").Value = 336
$ws.Range("F5").Value = 336
$ws.Range("G5").Value = 340
$ws.Range("H5").Value = 305
$ws.Range("I5").Value = 305
$ws.Range("J5").Value = $null
$ws.Range("K5").Value = 2551
$ws.Range("L5").Value = 1270
$ws.Range("M5").Value = 1281
$ws.Range("N5").Value = 1281
$ws.Range("O5").Value = $null
$ws.Range("P5").Value = 55
$ws.Range("Q5").Value = 108
$ws.Range("R5").Value = -110
$ws.Range("S5").Value = -124
$ws.Range("T5").Value = 117
$ws.Range("U5").Value = -9
$ws.Range("V5").Value = 668
$ws.Range("W5").Value = 10.36
$ws.Range("X5").Value = 9.41
$ws.Range("Y5").Value = 26.71
$ws.Range("Z5").Value = 12.64
$ws.Range("AA5").Value = 99.09
$ws.Range("AB5").Value = 2230.44
$ws.Range("AC5").Value = 2765
$ws.Range("AD5").Value = 6.87
$ws.Range("AE5").Value = 11612
$ws.Range("AF5").Value = 1.64
$ws.Range("AG5").Value = 180
$ws.Range("AH5").Value = 0.95
$ws.Range("AI5").Value = 6.51
$ws.Range("AJ5").Value = 11041708
# Row 6
$ws.Range("D6").Value = 3857
$ws.Range("E6").Value = 358
$ws.Range("F6").Value = 358
$ws.Range("G6").Value = 341
$ws.Range("H6").Value = 255
$ws.Range("I6").Value = 255
$ws.Range("K6").Value = 2585
$ws.Range("L6").Value = 1084
$ws.Range("M6").Value = 1502
$ws.Range("N6").Value = 1502
$ws.Range("P6").Value = 55
$ws.Range("Q6").Value = 235
$ws.Range("R6").Value = -129
$ws.Range("S6").Value = -211
$ws.Range("T6").Value = 128
$ws.Range("U6").Value = 106
$ws.Range("V6").Value = 470
$ws.Range("W6").Value = 9.279999999999999
$ws.Range("X6").Value = 6.61
$ws.Range("Y6").Value = 18.32
$ws.Range("Z6").Value = 9.93
$ws.Range("AA6").Value = 72.15000000000001
$ws.Range("AB6").Value = 2640.76
$ws.Range("AC6").Value = 2310
$ws.Range("AD6").Value = 6.15
$ws.Range("AE6").Value = 13609
$ws.Range("AF6").Value = 1.04
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 1.41
$ws.Range("AI6").Value = 8.65
$ws.Range("AJ6").Value = 11041708
# Row 7
$ws.Range("D7").Value = 3568
$ws.Range("E7").Value = 338
$ws.Range("G7").Value = 330
$ws.Range("H7").Value = 240
$ws.Range("I7").Value = 240
$ws.Range("K7").Value = 2830
$ws.Range("L7").Value = 1130
$ws.Range("M7").Value = 1700
$ws.Range("N7").Value = 1700
$ws.Range("P7").Value = 60
$ws.Range("Q7").Value = 400
$ws.Range("R7").Value = -150
$ws.Range("S7").Value = -20
$ws.Range("T7").Value = $null
$ws.Range("U7").Value = 180
$ws.Range("W7").Value = 9.470000000000001
$ws.Range("X7").Value = 6.73
$ws.Range("Y7").Value = 14.99
$ws.Range("Z7").Value = 8.859999999999999
$ws.Range("AA7").Value = 66.47
$ws.Range("AC7").Value = 2174
$ws.Range("AD7").Value = 4.69
$ws.Range("AE7").Value = 15404
$ws.Range("AF7").Value = 0.66
$ws.Range("AG7").Value = 200
$ws.Range("AH7").Value = 1.96
$ws.Range("AI7").Value = 9.199999999999999
# Row 8
$ws.Range("D8").Value = 3682
$ws.Range("E8").Value = 334
$ws.Range("G8").Value = 320
$ws.Range("H8").Value = 250
$ws.Range("I8").Value = 250
$ws.Range("K8").Value = 3030
$ws.Range("L8").Value = 1140
$ws.Range("M8").Value = 1890
$ws.Range("N8").Value = 1890
$ws.Range("P8").Value = 60
$ws.Range("Q8").Value = 310
$ws.Range("R8").Value = -160
$ws.Range("S8").Value = -10
$ws.Range("T8").Value = $null
$ws.Range("U8").Value = 190
$ws.Range("W8").Value = 9.07
$ws.Range("X8").Value = 6.79
$ws.Range("Y8").Value = 13.93
$ws.Range("Z8").Value = 8.529999999999999
$ws.Range("AA8").Value = 60.32
$ws.Range("AC8").Value = 2264
$ws.Range("AD8").Value = 4.51
$ws.Range("AE8").Value = 17126
$ws.Range("AF8").Value = 0.6
$ws.Range("AG8").Value = 200
$ws.Range("AH8").Value = 1.96
$ws.Range("AI8").Value = 8.83
# Row 9
$ws.Range("D9").Value = 3840
$ws.Range("E9").Value = 370
$ws.Range("G9").Value = 360
$ws.Range("H9").Value = 280
$ws.Range("I9").Value = 280
$ws.Range("K9").Value = 3240
$ws.Range("L9").Value = 1180
$ws.Range("M9").Value = 2060
$ws.Range("N9").Value = 2060
$ws.Range("P9").Value = 60
$ws.Range("Q9").Value = 320
$ws.Range("R9").Value = -160
$ws.Range("S9").Value = -10
$ws.Range("T9").Value = $null
$ws.Range("U9").Value = 220
$ws.Range("W9").Value = 9.630000000000001
$ws.Range("X9").Value = 7.29
$ws.Range("Y9").Value = 14.18
$ws.Range("Z9").Value = 8.93
$ws.Range("AA9").Value = 57.28
$ws.Range("AC9").Value = 2536
$ws.Range("AD9").Value = 4.02
$ws.Range("AE9").Value = 18666
$ws.Range("AF9").Value = 0.55
$ws.Range("AG9").Value = 200
$ws.Range("AH9").Value = 1.96
$ws.Range("AI9").Value = 7.89
